$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5:E6").Value = "ERROR"
